# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
# Every player row gets the team's season record: 69 wins, 93 losses, 0 ties.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): AD1="Wins", AE1="Losses", AF1="Ties" ---
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the existing header formatting (bold, thin border, centered / top-aligned)
$headerRange = $ws.Range("AD1:AF1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108  # xlCenter
$headerRange.VerticalAlignment = -4160    # xlTop
$headerRange.Borders.LineStyle = 1

# --- Data rows (2-48): season record repeated for every player ---
$wins = 69
$losses = 93
$ties = 0

for ($row = 2; $row -le 48; $row++) {
    $ws.Cells.Item($row, 30).Value = $wins    # column AD
    $ws.Cells.Item($row, 31).Value = $losses  # column AE
    $ws.Cells.Item($row, 32).Value = $ties    # column AF
}

Write-Output "Season record columns added."
